$wb = $excel.ActiveWorkbook

# --- KnightEntities: rebalance the Damage column (C) ---
$ws = $wb.Worksheets.Item("KnightEntities")

$ws.Range("C2:C6").Value   = 3
$ws.Range("C7:C11").Value  = 3.6
$ws.Range("C12:C16").Value = 4.3
$ws.Range("C17:C21").Value = 5

# --- Switch the active/selected tab from TowerEntities to KnightEntities,
#     and move the cell selection on KnightEntities to M18 ---
$ws.Activate()
$ws.Range("M18").Select()
